$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 174-176 (monthly revisions) ---
$ws.Cells.Item(174, 9).Value  = 4089   # I174
$ws.Cells.Item(174, 10).Value = 2114   # J174

$ws.Cells.Item(175, 8).Value  = 3709   # H175
$ws.Cells.Item(175, 9).Value  = 2218   # I175
$ws.Cells.Item(175, 10).Value = 1989   # J175
$ws.Cells.Item(175, 11).Value = 946    # K175
$ws.Cells.Item(175, 12).Value = 0      # L175

$ws.Cells.Item(176, 7).Value  = 3668   # G176
$ws.Cells.Item(176, 8).Value  = 2372   # H176
$ws.Cells.Item(176, 9).Value  = 2373   # I176
$ws.Cells.Item(176, 10).Value = 1439   # J176
$ws.Cells.Item(176, 11).Value = 1227   # K176

# --- Append new row 177 for period 01-08-2021 ---
# Column A holds a date-like label that must stay plain text (matching
# the existing "dd-mm-yyyy" text labels used throughout column A).
# Writing it straight to the destination cell makes Excel's smart-entry
# auto-detect it as a date and convert it to a serial number, so we
# stage the literal text as a formula result in a scratch cell and then
# copy only the *value* over, which keeps it a genuine text value and
# does not touch styles.xml at all.
$scratch = $ws.Cells.Item(177, 20)
$scratch.Formula = '="01-08-2021"'
$scratch.Copy()
$ws.Cells.Item(177, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item(177, 2).Value  = 18300
$ws.Cells.Item(177, 3).Value  = 15822
$ws.Cells.Item(177, 4).Value  = 7193
$ws.Cells.Item(177, 5).Value  = 8629
$ws.Cells.Item(177, 6).Value  = 1443
$ws.Cells.Item(177, 7).Value  = 4728
$ws.Cells.Item(177, 8).Value  = 3724
$ws.Cells.Item(177, 9).Value  = 2900
$ws.Cells.Item(177, 10).Value = 2369
$ws.Cells.Item(177, 11).Value = 501
$ws.Cells.Item(177, 12).Value = 157
$ws.Cells.Item(177, 13).Value = 2477
$ws.Cells.Item(177, 14).Value = 1265
$ws.Cells.Item(177, 15).Value = 1212
